# "barcodes ozon worked version"
# The delivery template lost its extra blank data row: row 5 (the first of
# the two identical "box count" entry rows) is removed, so the totals row
# that used to sum D6:O6 / B5:B6 / C5:C6 / B20 now sums D5:O5 / B5:B5 /
# C5:C5 / B19 one row higher, and everything below shifts up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 5 and shift rows 6.. up by one. Excel COM automatically
# re-anchors the dependent SUM()/SUM()/SUM() formulas in the (now) row 6
# totals line and keeps the merged ranges in sync.
$ws.Rows(5).Delete()

# The author's active selection ended up on the new totals row (row 6,
# selected as a whole row) after the edit.
[void]$ws.Range("A6:XFD6").Select()
